$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P3").Value = "maa://21249 (94.67), maa://26254 (96.67), **maa://22738 (50.0)"
$ws.Range("X3").Value = "maa://27396 (84.36), maa://27484 (96.69), maa://27480 (83.33)"
$ws.Range("T4").Value = "maa://32509 (95.8), maa://27295 (86.84), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.57)"
$ws.Range("D5").Value = "maa://21245 (84.55), maa://22744 (81.48)"
$ws.Range("A8").Value = "更新日期：2025.04.05 17:12:18"
$ws.Range("D10").Value = "***maa://25695 (18.32), ***maa://39951 (13.33), ***maa://34206 (22.22), ***maa://39243 (25.0), *maa://45271 (53.49)"
$ws.Range("X11").Value = "maa://36713 (97.82)"
$ws.Range("H12").Value = "maa://21867 (90.06), **maa://45826 (33.33)"
$ws.Range("D13").Value = "maa://24999 (92.2), maa://36673 (92.41), maa://25001 (85.92)"
$ws.Range("H13").Value = "*maa://21248 (73.71), **maa://22728 (46.67)"
$ws.Range("H15").Value = "maa://24304 (87.95), maa://21478 (89.74)"
$ws.Range("AF15").Value = "maa://21364 (80.29), *maa://36666 (77.12), *maa://22766 (68.33)"
$ws.Range("X16").Value = "maa://28501 (98.11), maa://28051 (96.0)"
$ws.Range("T17").Value = "**maa://42324 (50.0)"
$ws.Range("D18").Value = "maa://24570 (97.01)"
$ws.Range("AB21").Value = "maa://21443 (81.61), ***maa://23820 (30.0)"
$ws.Range("L23").Value = "maa://39756 (95.69), maa://39875 (94.59)"
$ws.Range("AF23").Value = "maa://31489 (94.74)"
$ws.Range("D24").Value = "*maa://24368 (78.66), *maa://46650 (66.67)"
$ws.Range("T27").Value = "*maa://30624 (75.81)"
$ws.Range("X28").Value = "maa://39929 (90.69), maa://41749 (91.75), ***maa://39723 (13.89)"
$ws.Range("AF28").Value = "maa://36660 (92.6), *maa://36701 (66.67)"
$ws.Range("AB30").Value = "maa://42979 (97.16), maa://45822 (100.0), *maa://45045 (80.0)"
$ws.Range("P34").Value = "maa://48817 (89.47)"
$ws.Range("P38").Value = "*maa://24383 (69.23)"
$ws.Range("T39").Value = "*maa://45788 (80.0), maa://47079 (93.55), *maa://45790 (75.0)"
$ws.Range("H46").Value = "maa://35931 (92.0), maa://43901 (93.75)"
$ws.Range("H53").Value = "maa://32534 (94.29), **maa://32434 (33.33)"
$ws.Range("H55").Value = "maa://32532 (92.19)"
$ws.Range("H64").Value = "maa://44405 (86.21)"
